$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values
$ws.Range("B2").Value = 50
$ws.Range("D2").Value = 75

$ws.Range("B6").Value = 15
$ws.Range("B7").Value = 15
$ws.Range("B8").Value = 15

# Update the active selection on the sheet
$ws.Range("F6").Select()
